$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.443826
$ws.Range("H2").Value = 64.331478
$ws.Range("I2").Value = 0.6062978927103765
$ws.Range("J2").Value = 0.6062978927103765
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.394935666666666
$ws.Range("N2").Value = 4.184806999999999
$ws.Range("O2").Value = 0.638793578492805
$ws.Range("P2").Value = 0.6387935784928049
$ws.Range("Q2").Value = 29.91275771719399
$ws.Range("R2").Value = 269.214819454746
$ws.Range("S2").Value = 0.3872992005171082
$ws.Range("T2").Value = 0.3872992005171081
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.443826
$ws.Range("H3").Value = 64.331478
$ws.Range("I3").Value = 0.6062978927103765
$ws.Range("J3").Value = 0.6062978927103765
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7887676666666668
$ws.Range("N3").Value = 2.366303
$ws.Range("O3").Value = 0.3612064215071951
$ws.Range("P3").Value = 0.361206421507195
$ws.Range("Q3").Value = 16.914196598426
$ws.Range("R3").Value = 152.227769385834
$ws.Range("S3").Value = 0.2189986921932684
$ws.Range("T3").Value = 0.2189986921932683
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.641794
$ws.Range("H4").Value = 10.925382
$ws.Range("I4").Value = 0.1029672609675761
$ws.Range("J4").Value = 0.1029672609675761
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.394935666666666
$ws.Range("N4").Value = 4.184806999999999
$ws.Range("O4").Value = 0.638793578492805
$ws.Range("P4").Value = 0.6387935784928049
$ws.Range("Q4").Value = 5.080068341252666
$ws.Range("R4").Value = 45.720615071274
$ws.Range("S4").Value = 0.06577482510108044
$ws.Range("T4").Value = 0.06577482510108043
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.641794
$ws.Range("H5").Value = 10.925382
$ws.Range("I5").Value = 0.1029672609675761
$ws.Range("J5").Value = 0.1029672609675761
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7887676666666668
$ws.Range("N5").Value = 2.366303
$ws.Range("O5").Value = 0.3612064215071951
$ws.Range("P5").Value = 0.361206421507195
$ws.Range("Q5").Value = 2.872529355860667
$ws.Range("R5").Value = 25.85276420274601
$ws.Range("S5").Value = 0.03719243586649564
$ws.Range("T5").Value = 0.03719243586649563
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.28284533333333
$ws.Range("H6").Value = 30.848536
$ws.Range("I6").Value = 0.2907348463220475
$ws.Range("J6").Value = 0.2907348463220475
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.394935666666666
$ws.Range("N6").Value = 4.184806999999999
$ws.Range("O6").Value = 0.638793578492805
$ws.Range("P6").Value = 0.6387935784928049
$ws.Range("Q6").Value = 14.34390771028355
$ws.Range("R6").Value = 129.095169392552
$ws.Range("S6").Value = 0.1857195528746165
$ws.Range("T6").Value = 0.1857195528746165
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.28284533333333
$ws.Range("H7").Value = 30.848536
$ws.Range("I7").Value = 0.2907348463220475
$ws.Range("J7").Value = 0.2907348463220475
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.7887676666666668
$ws.Range("N7").Value = 2.366303
$ws.Range("O7").Value = 0.3612064215071951
$ws.Range("P7").Value = 0.361206421507195
$ws.Range("Q7").Value = 8.110775920267557
$ws.Range("R7").Value = 72.99698328240801
$ws.Range("S7").Value = 0.1050152934474311
$ws.Range("T7").Value = 0.1050152934474311
